# Add a new "2022-Q3" sheet right after "总计" (total) and before "2022-Q2".
$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# Header row for the new sheet (same layout/format as the other quarterly sheets).
$headerRange = $q3.Range("B1:H1")
Set-HeaderStyle $headerRange
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Fund holdings data for 2022-Q3 (index, code, name, scale, position, pct, value, rank).
$q3Data = @(
    @(0,  "000362", "国泰聚信价值优势灵活配置混合A", "27.52", "89.04", "3.80", "1.0458", 8),
    @(1,  "000363", "国泰聚信价值优势灵活配置混合C", "13.07", "89.04", "3.80", "0.4967", 8),
    @(2,  "020010", "国泰金牛创新混合",               "13.26", "86.27", "2.74", "0.3633", 10),
    @(3,  "012173", "国泰兴泽优选一年持有期混合A",     "8.41",  "88.23", "4.18", "0.3515", 4),
    @(4,  "012174", "国泰兴泽优选一年持有期混合C",     "6.17",  "88.23", "4.18", "0.2579", 4),
    @(5,  "200006", "长城消费增值混合",               "5.46",  "90.90", "4.69", "0.2561", 3),
    @(6,  "007835", "国泰鑫睿混合",                   "8.30",  "79.49", "2.98", "0.2473", 10),
    @(7,  "013890", "国泰睿毅三年持有期混合A",         "4.86",  "89.26", "4.17", "0.2027", 5),
    @(8,  "003516", "国泰融安多策略灵活配置混合A",     "8.00",  "90.01", "2.11", "0.1688", 10),
    @(9,  "005244", "国泰聚优价值灵活配置混合A",       "4.61",  "87.30", "3.01", "0.1388", 9),
    @(10, "002938", "中银证券健康产业灵活配置混合",     "1.98",  "92.72", "4.89", "0.0968", 4),
    @(11, "020023", "国泰事件驱动策略混合A",           "2.59",  "92.18", "2.34", "0.0606", 6),
    @(12, "005245", "国泰聚优价值灵活配置混合C",       "1.80",  "87.30", "3.01", "0.0542", 9),
    @(13, "008619", "永赢医药健康股票C",               "0.40",  "94.40", "5.07", "0.0203", 7),
    @(14, "013891", "国泰睿毅三年持有期混合C",         "0.45",  "89.26", "4.17", "0.0188", 5),
    @(15, "008618", "永赢医药健康股票A",               "0.24",  "94.40", "5.07", "0.0122", 7),
    @(16, "014960", "国泰融安多策略灵活配置混合C",     "0.15",  "90.01", "2.11", "0.0032", 10),
    @(17, "015592", "国泰事件驱动策略混合C",           "0.02",  "92.18", "2.34", "0.0005", 6)
)

$row = 2
foreach ($rec in $q3Data) {
    $aCell = $q3.Cells.Item($row, 1)
    Set-HeaderStyle $aCell
    $aCell.Value = $rec[0]

    # Columns B-G store values as text (codes with leading zeros, and
    # percentages/amounts formatted with a fixed number of decimals) -
    # a leading apostrophe forces text storage exactly like typing it
    # into Excel would.
    $q3.Cells.Item($row, 2).Value = "'" + $rec[1]
    $q3.Cells.Item($row, 3).Value = $rec[2]
    $q3.Cells.Item($row, 4).Value = "'" + $rec[3]
    $q3.Cells.Item($row, 5).Value = "'" + $rec[4]
    $q3.Cells.Item($row, 6).Value = "'" + $rec[5]
    $q3.Cells.Item($row, 7).Value = "'" + $rec[6]
    $q3.Cells.Item($row, 8).Value = $rec[7]
    $row = $row + 1
}

# Update the "总计" (total) summary sheet: prepend a 2022-Q3 row and shift
# the previously existing rows down by one.
$total = $wb.Worksheets.Item(1)

$total.Cells.Item(5, 1).Value = 3
Set-HeaderStyle $total.Cells.Item(5, 1)
$total.Cells.Item(5, 2).Value = "2021-Q2"
$total.Cells.Item(5, 3).Value = 2
$total.Cells.Item(5, 4).Value = 0.03

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q1"
$total.Cells.Item(4, 3).Value = 5
$total.Cells.Item(4, 4).Value = 0.62

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 6
$total.Cells.Item(3, 4).Value = 0.98

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 18
$total.Cells.Item(2, 4).Value = 3.8
